$d = $word.ActiveDocument

# The heading paragraph looks like:
#   [run: mc:AlternateContent floating text box, zero-width in the text model]
#   [run: "Intellij Gotchas"]
# We need to insert a new run containing just "i" immediately before the
# "Intellij Gotchas" run (i.e. right after the AlternateContent run), so the
# paragraph reads "iIntellij Gotchas".

# First locate the target text and remember where it starts.
$findRange = $d.Content
$findRange.Find.Execute("Intellij Gotchas", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $findRange.Start

# Replace the text in-place (rather than inserting into a collapsed, empty
# range at that offset) so the new character lands inside/after the existing
# run sequence -- i.e. after the floating drawing's run -- instead of being
# pushed in front of it.
$d.Content.Find.Execute("Intellij Gotchas", $false, $false, $false, $false, `
    $false, $true, 1, $false, "iIntellij Gotchas", 2) | Out-Null

# The newly inserted "i" now sits at $startPos. Give it its own run by
# explicitly (re)applying the same character formatting used by the rest of
# the heading: Arial 14pt bold (incl. complex-script variants), no proofing,
# and the en-GB far-east language tag.
$iRange = $d.Range($startPos, $startPos + 1)
$iRange.Font.Name = "Arial"
$iRange.Font.NameBi = "Arial"
$iRange.Font.Bold = $true
$iRange.Font.BoldBi = $true
$iRange.Font.Size = 14
$iRange.NoProofing = 1
$iRange.LanguageIDFarEast = "en-GB"
